# Applies the change described by the diff: inserts a new data row at
# row 426 (pushing the existing rows 426:473 down to 427:474) and fills
# the newly inserted row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 426; this shifts rows 426:473 down to 427:474
# and grows the sheet dimension from A1:R473 to A1:R474 automatically.
$ws.Rows(426).Insert()

# Populate the newly inserted row 426 with its data.
$ws.Range("A426").Value = 6
$ws.Range("B426").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C426").Value = "Metropolitana"
$ws.Range("D426").Value = 44449
$ws.Range("E426").Value = 13
$ws.Range("F426").Value = 100112031
$ws.Range("G426").Value = "Poroto verde"
$ws.Range("H426").Value = "Magnum"
$ws.Range("I426").Value = "Primera"
$ws.Range("J426").Value = 400
$ws.Range("K426").Value = 37000
$ws.Range("L426").Value = 40000
$ws.Range("M426").Value = 38725
$ws.Range("N426").Value = "$/malla 25 kilos"
$ws.Range("O426").Value = "Perú"
$ws.Range("P426").Value = 1549
$ws.Range("Q426").Value = 25
$ws.Range("R426").Value = "Hortaliza"
